# Release version 0.3.5 (#83)
# Update the Reqs.xlsx test fixture: the "Type" column (C) values change
# from "SW" to "MS" for every requirement row (rows 1-20 on Sheet1),
# including row 12 which previously had no value in column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 1; $r -le 20; $r++) {
    $ws.Cells.Item($r, 3).Value = "MS"
}

# Move the active selection from D1:D20 to I12, matching the saved
# worksheet view state in the updated workbook.
$ws.Range("I12").Select()
